$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "27.158.01"
$ws.Range("E2").Value = "  +0.21%  "
$ws.Range("D3").Value = "1.904.99"
$ws.Range("E3").Value = "  +0.64%  "
$ws.Range("E4").Value = "  +0.22%  "
$ws.Range("D5").Value = "'306.63"
$ws.Range("D5").Style = "Normal"
$ws.Range("E6").Value = "  +0.22%  "
$ws.Range("D7").Value = "'0.5236"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  +1.84%  "
$ws.Range("D8").Value = "'0.3762"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  +0.24%  "
$ws.Range("D9").Value = "'0.07262"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  +0.82%  "
$ws.Range("D10").Value = "'21.19"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  -0.14%  "
$ws.Range("D11").Value = "'0.9067"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  +0.18%  "
$ws.Range("D12").Value = "'0.08485"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  +11.11%  "
$ws.Range("D13").Value = "1.919.14"
$ws.Range("E13").Value = "  +1.34%  "
$ws.Range("D14").Value = "'96.97"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  +2.07%  "
$ws.Range("D15").Value = "'5.294"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  +0.57%  "
$ws.Range("D16").Value = "'1.000"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  +0.15%  "
$ws.Range("D17").Value = "'0.000008701"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  +2.49%  "
$ws.Range("D18").Value = "'14.55"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  +0.75%  "
$ws.Range("E19").Value = "  +0.23%  "
$ws.Range("D20").Value = "27.194.26"
$ws.Range("E20").Value = "  +0.25%  "
$ws.Range("D21").Value = "'5.085"
$ws.Range("D21").Style = "Normal"
$ws.Range("D22").Value = "2.155.23"
$ws.Range("E22").Value = "  +0.94%  "
$ws.Range("D23").Value = "'10.64"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  +0.77%  "
$ws.Range("D24").Value = "'6.440"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  +0.52%  "
$ws.Range("D25").Value = "'2.320"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +1.61%  "
$ws.Range("D26").Value = "'146.88"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  +0.90%  "
$ws.Range("D27").Value = "'18.26"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  +1.21%  "
$ws.Range("D28").Value = "'1.744"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  -1.32%  "
$ws.Range("D29").Value = "'115.15"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  +0.59%  "
$ws.Range("B30").Value = "InternetComputer(DFINITY)"
$ws.Range("C30").Value = "https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp"
$ws.Range("D30").Value = "'4.828"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  -0.02%  "
$ws.Range("B31").Value = "Filecoin"
$ws.Range("C31").Value = "https://coinranking.com/coin/ymQub4fuB+filecoin-fil"
$ws.Range("D31").Value = "'4.920"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  -0.53%  "
$ws.Range("D32").Value = "'0.09306"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  +1.37%  "
$ws.Range("E33").Value = "  +2.88%  "
$ws.Range("D34").Value = "'0.05058"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  -0.55%  "
$ws.Range("E35").Value = "  +0.43%  "
$ws.Range("D36").Value = "'3.440"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  +4.55%  "
$ws.Range("D37").Value = "'2.945"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  -1.42%  "
$ws.Range("D38").Value = "'2.600"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  -0.50%  "
$ws.Range("D39").Value = "'0.5714"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  +2.03%  "
$ws.Range("D40").Value = "'0.02003"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  +0.25%  "
$ws.Range("D41").Value = "'1.075"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  -0.11%  "
$ws.Range("D42").Value = "'9.126"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  +0.28%  "
$ws.Range("D43").Value = "'6.632"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  -0.28%  "
$ws.Range("D44").Value = "'115.96"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  -1.47%  "
$ws.Range("D45").Value = "'0.1515"
$ws.Range("D45").Style = "Normal"
$ws.Range("D46").Value = "'0.4864"
$ws.Range("D46").Style = "Normal"
$ws.Range("B47").Value = "PaxDollar"
$ws.Range("C47").Value = "https://coinranking.com/coin/JCKLgWPAF+paxdollar-usdp"
$ws.Range("D47").Value = "'1.0000"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  +0.25%  "
$ws.Range("B48").Value = "EnergySwap"
$ws.Range("C48").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Range("D48").Value = "'10.15"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  -0.38%  "
$ws.Range("D49").Value = "'1.623"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  +1.52%  "
$ws.Range("D50").Value = "'37.75"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  +0.65%  "
$ws.Range("D51").Value = "'64.22"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  +0.31%  "
